$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in A1 (was "contenedor", now "mobile version")
$ws.Range("A1").Value = "mobile version"

# Add the new "desktop" section
$ws.Range("A11").Value = "desktop"
$ws.Range("B12").Value = 870
$ws.Range("B13").Value = 1805
$ws.Range("B14").Value = 1195
$ws.Range("B15").Value = 686
$ws.Range("B16").Formula = "=SUM(B12:B15)"

# Update selection to match the new active cell
$ws.Range("B16").Select()
